# Applies the StructureDefinition-process-name.xlsx update:
#  - Metadata sheet: Version, Date, Publisher, Jurisdiction fields updated;
#    the duplicate "Contact" row is removed.
#  - Elements sheet: row 2 (root Extension) Short/Definition text updated.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# remove it so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Process Name"
$elements.Range("L2").Value = "The name of the process or service that has produced the data held in the FHIR resource or element. Example: If an IBM analytic service has been the producer, use the process-name as the name of that analytic service."
